# Updating FoodLog 15/04/2018 file
# Appends rows 99-105 (09-04-2018 .. 15-04-2018) to the daily food log on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -----------------------------------------------------
# columns: Date, Calories In, Fat(g), Fiber(g), Carbs(g), Sodium(mg), Protein(g), Water(ml)
$newRows = @(
    @{ Row = 99;  Date = "09-04-2018"; Cal = 1256; Fat = 63; Fiber = 2;  Carbs = 101; Sodium = 2181; Protein = 66;  Water = 2500 },
    @{ Row = 100; Date = "10-04-2018"; Cal = 1299; Fat = 31; Fiber = 8;  Carbs = 134; Sodium = 1855; Protein = 110; Water = 2000 },
    @{ Row = 101; Date = "11-04-2018"; Cal = 1545; Fat = 62; Fiber = 9;  Carbs = 10;  Sodium = 2225; Protein = 130; Water = 1500 },
    @{ Row = 102; Date = "12-04-2018"; Cal = 1548; Fat = 54; Fiber = 3;  Carbs = 117; Sodium = 2621; Protein = 130; Water = 1000 },
    @{ Row = 103; Date = "13-04-2018"; Cal = 1356; Fat = 37; Fiber = 13; Carbs = 168; Sodium = 1771; Protein = 81;  Water = 1750 },
    @{ Row = 104; Date = "14-04-2018"; Cal = 1264; Fat = 31; Fiber = 12; Carbs = 182; Sodium = 1709; Protein = 55;  Water = 2000 },
    @{ Row = 105; Date = "15-04-2018"; Cal = 1209; Fat = 31; Fiber = 16; Carbs = 166; Sodium = 1511; Protein = 62;  Water = 3000 }
)

$firstRow = 99
$lastRow  = 105

# Column A (dates) must stay plain text, matching the existing log entries,
# instead of being auto-converted into date serials - format as Text first.
$ws.Range("A$firstRow`:A$lastRow").NumberFormat = "@"

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Date
    $ws.Range("B$row").Value = $r.Cal
    $ws.Range("C$row").Value = $r.Fat
    $ws.Range("D$row").Value = $r.Fiber
    $ws.Range("E$row").Value = $r.Carbs
    $ws.Range("F$row").Value = $r.Sodium
    $ws.Range("G$row").Value = $r.Protein
    $ws.Range("H$row").Value = $r.Water
}

# Match formatting used by the rest of the log.
$ws.Range("A$firstRow`:A$lastRow").Style = "Normal 2"
$ws.Range("B$firstRow`:B$lastRow").Style = "Normal 2"
$ws.Range("B$firstRow`:B$lastRow").NumberFormat = "0"
$ws.Range("H$firstRow`:H$lastRow").NumberFormat = "0"

# Extend the WaterTargetAchieved / UnderEaten formulas down through the new rows.
foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("I$row").Formula = "=IF(H$row>=2200,""Yes"",""No"")"
    $ws.Range("J$row").Formula = "=IF(B$row<=1800,""Yes"",""No"")"
}

# Move the selection to reflect where the user ended up after typing the new data.
$ws.Range("H106").Select()

$wb.Application.Calculate()
